$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.302.53'
$ws.Range("E2").Value = '  -3.43%  '

# Row 3
$ws.Range("D3").Value = '1.811.92'
$ws.Range("E3").Value = '  -3.73%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.77'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.92%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4215'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.54%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3563'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07156'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.85%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8477'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.78%  '

# Row 11
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.19'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.94%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.902.43'
$ws.Range("E12").Value = '  -1.35%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.340'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.78%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.383'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -4.08%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06923'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.81%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.18%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.37'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.03%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008827'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.70%  '

# Row 19
$ws.Range("E19").Value = '  +0.08%  '

# Row 20
$ws.Range("E20").Value = '  -3.28%  '

# Row 21
$ws.Range("D21").Value = '27.423.27'
$ws.Range("E21").Value = '  -3.40%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.095'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.27%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.93'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.20%  '

# Row 24
$ws.Range("D24").Value = '2.068.77'
$ws.Range("E24").Value = '  -4.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.966'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.30%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.50'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.58%  '

# Row 27
$ws.Range("E27").Value = '  -2.96%  '

# Row 28
$ws.Range("E28").Value = '  -6.76%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.38'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.53%  '

# Row 30
$ws.Range("E30").Value = '  -9.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08910'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.96%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7425'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.52%  '

# Row 33
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.937'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.12%  '

# Row 34
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.484'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.00%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.111'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.21%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.003'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.10%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.073'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.59%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05215'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.01%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01908'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.99%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.773'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.99%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1643'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.65%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4995'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.76%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.306'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -8.85%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.219'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.00%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.34'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.27%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '105.23'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.84%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06412'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.10%  '

# Row 48
$ws.Range("E48").Value = '  +0.14%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4601'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.92%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.604'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.87%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.26'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.45%  '
